$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1539
$ws.Range("J43").Value = 1600
$ws.Range("L43").Value = 1600
$ws.Range("N43").Value = -1738
$ws.Range("H76").Value = 4691403.5
$ws.Range("I76").Value = 11715712
$ws.Range("J76").Value = 8531.666999999999
$ws.Range("K76").Value = 11715712
$ws.Range("L76").Value = 8531.666999999999
$ws.Range("M76").Value = -11715397
$ws.Range("N76").Value = -9161.666999999999
$ws.Range("H79").Value = 4691403.5
$ws.Range("I79").Value = 11715712
$ws.Range("J79").Value = 8531.666999999999
$ws.Range("K79").Value = 11715712
$ws.Range("L79").Value = 8531.666999999999
$ws.Range("M79").Value = -11714620
$ws.Range("N79").Value = -10715.667
$ws.Range("H137").Value = 34776.7
$ws.Range("J137").Value = 112676.89
$ws.Range("L137").Value = 338030.67
$ws.Range("N137").Value = -343130.67
$ws.Range("H141").Value = 905876.9399999999
$ws.Range("J141").Value = 5662.857
$ws.Range("L141").Value = 16988.571
$ws.Range("N141").Value = -27348.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3323344.2
$ws.Range("I2").Value = 4652082
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 4652082
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -4651969
$ws.Range("N2").Value = -1726
$ws.Range("H45").Value = 1568.8235
$ws.Range("I45").Value = 1052.1666
$ws.Range("J45").Value = 1850.6364
$ws.Range("K45").Value = 1052.1666
$ws.Range("L45").Value = 1850.6364
$ws.Range("M45").Value = -675.1666
$ws.Range("N45").Value = -2604.6364
$ws.Range("H61").Value = 2209.1482
$ws.Range("I61").Value = 1170.7142
$ws.Range("K61").Value = 1170.7142
$ws.Range("M61").Value = -958.7141999999999
$ws.Range("H74").Value = 1175.3829
$ws.Range("I74").Value = 989.8919
$ws.Range("K74").Value = 989.8919
$ws.Range("M74").Value = -115.8919
$ws.Range("H77").Value = 1175.3829
$ws.Range("I77").Value = 989.8919
$ws.Range("K77").Value = 4949.4595
$ws.Range("M77").Value = -581.4594999999999
$ws.Range("H116").Value = 3323344.2
$ws.Range("I116").Value = 4652082
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 4652082
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = -4649788
$ws.Range("N116").Value = -6088
$ws.Range("H132").Value = 1394.8158
$ws.Range("I132").Value = 932.6429000000001
$ws.Range("K132").Value = 2797.9287
$ws.Range("M132").Value = -267.9287000000004
$ws.Range("H136").Value = 2209.1482
$ws.Range("I136").Value = 1170.7142
$ws.Range("K136").Value = 3512.1426
$ws.Range("M136").Value = -962.1425999999997
$ws.Range("H139").Value = 30084
$ws.Range("J139").Value = 30084
$ws.Range("L139").Value = 30084
$ws.Range("N139").Value = -40364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3323344.2
$ws.Range("I3").Value = 4652082
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 4652082
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -4651968
$ws.Range("N3").Value = -1728
$ws.Range("H134").Value = 4671.641
$ws.Range("I134").Value = 5078.8237
$ws.Range("J134").Value = 1902.8
$ws.Range("K134").Value = 15236.4711
$ws.Range("L134").Value = 5708.4
$ws.Range("M134").Value = -12701.4711
$ws.Range("N134").Value = -10778.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2168.3076
$ws.Range("I31").Value = 1727.4286
$ws.Range("K31").Value = 1727.4286
$ws.Range("M31").Value = -1432.4286
$ws.Range("H34").Value = 2168.3076
$ws.Range("I34").Value = 1727.4286
$ws.Range("K34").Value = 1727.4286
$ws.Range("M34").Value = -1525.4286
$ws.Range("H35").Value = 200.33333
$ws.Range("I35").Value = 200.33333
$ws.Range("K35").Value = 200.33333
$ws.Range("M35").Value = 93.66667000000001
$ws.Range("H58").Value = 1451262.6
$ws.Range("I58").Value = 2175338
$ws.Range("K58").Value = 2175338
$ws.Range("M58").Value = -2175135
$ws.Range("H132").Value = 1753.24
$ws.Range("I132").Value = 1115.2778
$ws.Range("K132").Value = 3345.8334
$ws.Range("M132").Value = -815.8334000000004
$ws.Range("H134").Value = 1761.5172
$ws.Range("I134").Value = 1477.2222
$ws.Range("K134").Value = 4431.6666
$ws.Range("M134").Value = -1896.6666
$ws.Range("H136").Value = 1451262.6
$ws.Range("I136").Value = 2175338
$ws.Range("K136").Value = 6526014
$ws.Range("M136").Value = -6523464

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1787.625
$ws.Range("I118").Value = 1354.5
$ws.Range("J118").Value = 2220.75
$ws.Range("K118").Value = 4063.5
$ws.Range("L118").Value = 6662.25
$ws.Range("M118").Value = -2820.5
$ws.Range("N118").Value = -9148.25
$ws.Range("H122").Value = 1009.1818
$ws.Range("J122").Value = 1296.8334
$ws.Range("L122").Value = 11671.5006
$ws.Range("N122").Value = -16571.5006
$ws.Range("H132").Value = 1764.1428
$ws.Range("I132").Value = 1379.8
$ws.Range("K132").Value = 12418.2
$ws.Range("M132").Value = -9888.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5500
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -9831
$ws.Range("N10").Value = -1338
$ws.Range("H19").Value = 63143.25
$ws.Range("I19").Value = 55555
$ws.Range("J19").Value = 65672.664
$ws.Range("K19").Value = 55555
$ws.Range("L19").Value = 65672.664
$ws.Range("M19").Value = -55267
$ws.Range("N19").Value = -66248.664
$ws.Range("H132").Value = 1204414.1
$ws.Range("I132").Value = 1833452.4
$ws.Range("J132").Value = 3522.9092
$ws.Range("K132").Value = 5500357.199999999
$ws.Range("L132").Value = 10568.7276
$ws.Range("M132").Value = -5497827.199999999
$ws.Range("N132").Value = -15628.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1122.5
$ws.Range("I100").Value = 1345
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 1345
$ws.Range("L100").Value = 900
$ws.Range("M100").Value = -804
$ws.Range("N100").Value = -1982
$ws.Range("H132").Value = 1679.7273
$ws.Range("I132").Value = 1414.0476
$ws.Range("J132").Value = 1922.3043
$ws.Range("K132").Value = 4242.142800000001
$ws.Range("L132").Value = 5766.9129
$ws.Range("M132").Value = -1712.142800000001
$ws.Range("N132").Value = -10826.9129
$ws.Range("H136").Value = 2374.2593
$ws.Range("I136").Value = 1430.25
$ws.Range("J136").Value = 5071.4287
$ws.Range("K136").Value = 4290.75
$ws.Range("L136").Value = 15214.2861
$ws.Range("M136").Value = -1740.75
$ws.Range("N136").Value = -20314.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30452
$ws.Range("H132").Value = 1223.1794
$ws.Range("I132").Value = 927.5862
$ws.Range("J132").Value = 2080.4
$ws.Range("K132").Value = 2782.7586
$ws.Range("L132").Value = 6241.200000000001
$ws.Range("M132").Value = -252.7586000000001
$ws.Range("N132").Value = -11301.2
$ws.Range("H136").Value = 21370326
$ws.Range("I136").Value = 27780502
$ws.Range("K136").Value = 83341506
$ws.Range("M136").Value = -83338956
